$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.190906763076782
$ws.Range("B1").Value = 2.473937034606934
$ws.Range("D1").Value = 2.278348922729492
$ws.Range("E1").Value = 1.179724931716919
